$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 22
$ws.Range("E13").Value = "criado repositorys e daos de forma diferente ao que foi abordado em outro curso"
$ws.Range("D13").Value = "22. Incluindo DAO'S especificos"
$ws.Range("C13").Value = "4. Camada de Persistência"

$ws.Range("E13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 30

$ws.Range("B14").Select()
